# Automatic update of files.
# Update the "last updated" date column (C) for every data row from
# 2023-09-12 (serial 45181) to 2023-09-13 (serial 45182).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 72; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
